$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new content: cell A1 gets the text "lista 2"
$ws.Range("A1").Value = "lista 2"

# Match the cursor/selection state left in the sheet (H8 selected)
[void]$ws.Range("H8").Select()
